# Add "import ontologies from excel" metadata to the Metadata sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# B2: Ontology IRI
$ws.Range("B2").Value = "http://emmo.info/emmo/domain/onto"

# B3: Ontology version IRI
$ws.Range("B3").Value = "http://emmo.info/emmo/domain/onto/0.01"

# B4: Ontology version Info (numeric)
$ws.Range("B4").Value = 0.01

# B13: Imported ontologies (full IRI)
$ws.Range("B13").Value = "https://raw.githubusercontent.com/emmo-repo/emmo-repo.github.io/master/versions/1.0.0-beta/emmo-inferred-chemistry2.ttl"

# Column B widened to fit the long URL that was just typed in (bestFit-style autosize)
$ws.Columns.Item(2).ColumnWidth = 119.3

# Metadata sheet becomes the active tab, with A13 selected (matches the
# author's last interaction with the sheet after filling in the import row)
$ws.Activate()
$ws.Range("A13").Select()
